$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 7924.625
$ws.Range("I28").Value = 10660
$ws.Range("J28").Value = 3365.6667
$ws.Range("K28").Value = 10660
$ws.Range("L28").Value = 3365.6667
$ws.Range("M28").Value = -10175
$ws.Range("N28").Value = -4335.6667
# Row 51
$ws.Range("H51").Value = 9996.333000000001
$ws.Range("I51").Value = 9989
$ws.Range("K51").Value = 9989
$ws.Range("M51").Value = -9505
# Row 96
$ws.Range("H96").Value = 2120.4
$ws.Range("I96").Value = 2069.5
$ws.Range("J96").Value = 2324
$ws.Range("K96").Value = 6208.5
$ws.Range("L96").Value = 6972
$ws.Range("M96").Value = -4835.5
$ws.Range("N96").Value = -9718
# Row 98
$ws.Range("H98").Value = 38498.53
$ws.Range("I98").Value = 55749.9
$ws.Range("J98").Value = 13853.714
$ws.Range("K98").Value = 55749.9
$ws.Range("L98").Value = 13853.714
$ws.Range("M98").Value = -54251.9
$ws.Range("N98").Value = -16849.714
# Row 113
$ws.Range("H113").Value = 13294.786
$ws.Range("I113").Value = 20547.25
$ws.Range("K113").Value = 20547.25
$ws.Range("M113").Value = -17293.25
# Row 122
$ws.Range("H122").Value = 38498.53
$ws.Range("I122").Value = 55749.9
$ws.Range("J122").Value = 13853.714
$ws.Range("K122").Value = 167249.7
$ws.Range("L122").Value = 41561.142
$ws.Range("M122").Value = -164799.7
$ws.Range("N122").Value = -46461.142
# Row 130
$ws.Range("H130").Value = 70000
$ws.Range("J130").Value = 70000
$ws.Range("L130").Value = 70000
$ws.Range("N130").Value = -80040
# Row 137
$ws.Range("H137").Value = 629967.25
$ws.Range("I137").Value = 815265.3
$ws.Range("J137").Value = 27748.5
$ws.Range("K137").Value = 2445795.9
$ws.Range("L137").Value = 83245.5
$ws.Range("M137").Value = -2443245.9
$ws.Range("N137").Value = -88345.5
# Row 138
$ws.Range("H138").Value = 4438.264
$ws.Range("J138").Value = 5164.0425
$ws.Range("L138").Value = 15492.1275
$ws.Range("N138").Value = -25772.1275

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 7256
$ws.Range("I2").Value = 9082.532999999999
$ws.Range("K2").Value = 9082.532999999999
$ws.Range("M2").Value = -8969.532999999999
# Row 45
$ws.Range("H45").Value = 130191.766
$ws.Range("I45").Value = 197137.27
$ws.Range("K45").Value = 197137.27
$ws.Range("M45").Value = -196760.27
# Row 61
$ws.Range("H61").Value = 8761.1
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
# Row 74
$ws.Range("H74").Value = 3428.8928
$ws.Range("I74").Value = 1704.7826
$ws.Range("J74").Value = 11359.8
$ws.Range("K74").Value = 1704.7826
$ws.Range("L74").Value = 11359.8
$ws.Range("M74").Value = -830.7826
$ws.Range("N74").Value = -13107.8
# Row 77
$ws.Range("H77").Value = 3428.8928
$ws.Range("I77").Value = 1704.7826
$ws.Range("J77").Value = 11359.8
$ws.Range("K77").Value = 8523.913
$ws.Range("L77").Value = 56799
$ws.Range("M77").Value = -4155.913
$ws.Range("N77").Value = -65535
# Row 97
$ws.Range("H97").Value = 6536.2856
$ws.Range("I97").Value = 7559.7646
$ws.Range("K97").Value = 7559.7646
$ws.Range("M97").Value = -7063.7646
# Row 110
$ws.Range("H110").Value = 3400
$ws.Range("I110").Value = 3400
$ws.Range("K110").Value = 3400
$ws.Range("M110").Value = -1355
# Row 116
$ws.Range("H116").Value = 7256
$ws.Range("I116").Value = 9082.532999999999
$ws.Range("K116").Value = 9082.532999999999
$ws.Range("M116").Value = -6788.532999999999
# Row 136
$ws.Range("H136").Value = 8761.1
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 7256
$ws.Range("I3").Value = 9082.532999999999
$ws.Range("K3").Value = 9082.532999999999
$ws.Range("M3").Value = -8968.532999999999
# Row 20
$ws.Range("H20").Value = 2964.55
$ws.Range("I20").Value = 1749.5
$ws.Range("K20").Value = 1749.5
$ws.Range("M20").Value = -1502.5
# Row 94
$ws.Range("H94").Value = 2562.139
$ws.Range("I94").Value = 2073.32
$ws.Range("J94").Value = 3673.0908
$ws.Range("K94").Value = 2073.32
$ws.Range("L94").Value = 3673.0908
$ws.Range("M94").Value = -1622.32
$ws.Range("N94").Value = -4575.0908
# Row 134
$ws.Range("H134").Value = 14642.75
$ws.Range("I134").Value = 19233.5
$ws.Range("K134").Value = 57700.5
$ws.Range("M134").Value = -55165.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2581.7058
$ws.Range("I16").Value = 2799.6365
$ws.Range("K16").Value = 2799.6365
$ws.Range("M16").Value = -2512.6365
# Row 99
$ws.Range("H99").Value = 254932
$ws.Range("I99").Value = 421094.16
$ws.Range("K99").Value = 421094.16
$ws.Range("M99").Value = -419596.16
# Row 107
$ws.Range("H107").Value = 83362136
$ws.Range("J107").Value = 1999.6
$ws.Range("L107").Value = 1999.6
$ws.Range("N107").Value = -5839.6
# Row 113
$ws.Range("H113").Value = 2581.7058
$ws.Range("I113").Value = 2799.6365
$ws.Range("K113").Value = 2799.6365
$ws.Range("M113").Value = -629.6365000000001
# Row 126
$ws.Range("H126").Value = 254932
$ws.Range("I126").Value = 421094.16
$ws.Range("K126").Value = 1263282.48
$ws.Range("M126").Value = -1260812.48

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 8288286
$ws.Range("I4").Value = 8288286
$ws.Range("K4").Value = 24864858
$ws.Range("M4").Value = -24864746
# Row 46
$ws.Range("H46").Value = 2435.3635
$ws.Range("I46").Value = 498.33334
$ws.Range("J46").Value = 4759.8
$ws.Range("K46").Value = 1495.00002
$ws.Range("L46").Value = 14279.4
$ws.Range("M46").Value = -1404.00002
$ws.Range("N46").Value = -14461.4
# Row 48
$ws.Range("H48").Value = 3679.1428
$ws.Range("I48").Value = 1310.8
$ws.Range("K48").Value = 3932.4
$ws.Range("M48").Value = -3682.4
# Row 51
$ws.Range("H51").Value = 502166.94
$ws.Range("I51").Value = 770087.9399999999
$ws.Range("J51").Value = 4599.4287
$ws.Range("K51").Value = 2310263.82
$ws.Range("L51").Value = 13798.2861
$ws.Range("M51").Value = -2309803.82
$ws.Range("N51").Value = -14718.2861
# Row 130
$ws.Range("H130").Value = 4450
$ws.Range("I130").Value = 4450
$ws.Range("K130").Value = 13350
$ws.Range("M130").Value = -8330
# Row 131
$ws.Range("H131").Value = 100001736
$ws.Range("J131").Value = 1933
$ws.Range("L131").Value = 5799
$ws.Range("N131").Value = -15879
# Row 134
$ws.Range("H134").Value = 5332.615
$ws.Range("I134").Value = 1932.4
$ws.Range("J134").Value = 16666.666
$ws.Range("K134").Value = 5797.200000000001
$ws.Range("L134").Value = 49999.99800000001
$ws.Range("M134").Value = -727.2000000000007
$ws.Range("N134").Value = -60139.99800000001
# Row 140
$ws.Range("H140").Value = 13140.533
$ws.Range("I140").Value = 13140.533
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 39421.599
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -34241.599
$ws.Range("N140").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 12986.613
$ws.Range("I122").Value = 11206.6
$ws.Range("K122").Value = 33619.8
$ws.Range("M122").Value = -31169.8
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 4739.579
$ws.Range("I122").Value = 4467.4614
$ws.Range("J122").Value = 5329.1665
$ws.Range("K122").Value = 13402.3842
$ws.Range("L122").Value = 15987.4995
$ws.Range("M122").Value = -10952.3842
$ws.Range("N122").Value = -20887.4995

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 3199.2856
$ws.Range("J96").Value = 3711
$ws.Range("L96").Value = 3711
$ws.Range("N96").Value = -6457
# Row 100
$ws.Range("H100").Value = 25596.143
$ws.Range("I100").Value = 15038.066
$ws.Range("J100").Value = 51991.332
$ws.Range("K100").Value = 30076.132
$ws.Range("L100").Value = 103982.664
$ws.Range("M100").Value = -29535.132
$ws.Range("N100").Value = -105064.664
# Row 126
$ws.Range("H126").Value = 41799.363
$ws.Range("I126").Value = 59971.855
$ws.Range("K126").Value = 179915.565
$ws.Range("M126").Value = -177445.565
# Row 132
$ws.Range("H132").Value = 71995
$ws.Range("I132").Value = 157497.5
$ws.Range("J132").Value = 14993.333
$ws.Range("K132").Value = 472492.5
$ws.Range("L132").Value = 44979.999
$ws.Range("M132").Value = -469962.5
